# Daily attendance processing - 2026-01-06 23:02:34
# Swap the order of "Recorded By" entries in column G from
# "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}
